$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2-28
# from serial 45233 (2023-11-03) to serial 45243 (2023-11-13)
$ws.Range("C2:C28").Value = 45243
